$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -12.37
$ws.Range("D4").Value = -7.795999999999999
$ws.Range("C6").Value = -11.714
$ws.Range("C7").Value = -12.995
$ws.Range("C8").Value = -12.282
$ws.Range("D8").Value = -7.719000000000001
$ws.Range("D9").Value = -8.026999999999999
$ws.Range("D12").Value = -7.451000000000001
$ws.Range("C16").Value = -12.429
$ws.Range("D17").Value = -8.331999999999999
$ws.Range("D18").Value = -8.615
$ws.Range("D19").Value = -7.731999999999999
$ws.Range("C20").Value = -12.241
$ws.Range("D20").Value = -7.375
$ws.Range("C21").Value = -12.862
$ws.Range("D26").Value = -7.371
$ws.Range("C28").Value = -12.25
$ws.Range("C29").Value = -12.148
$ws.Range("C30").Value = -11.162
$ws.Range("D31").Value = -7.991
$ws.Range("C32").Value = -12.141
$ws.Range("D39").Value = -7.436999999999999
$ws.Range("C40").Value = -12.241
$ws.Range("D40").Value = -7.558
$ws.Range("D41").Value = -7.685
$ws.Range("D42").Value = -7.859
$ws.Range("D43").Value = -7.947
$ws.Range("C46").Value = -13.556
$ws.Range("D47").Value = -7.475
$ws.Range("D48").Value = -7.585000000000001
$ws.Range("C51").Value = -11.438
$ws.Range("C52").Value = -11.538
$ws.Range("D54").Value = -7.835000000000001
$ws.Range("C57").Value = -13.848
$ws.Range("C59").Value = -11.455
$ws.Range("C62").Value = -13.62
$ws.Range("D62").Value = -8.356999999999999
$ws.Range("D63").Value = -7.25
$ws.Range("D64").Value = -7.580999999999999
$ws.Range("C66").Value = -11.855
$ws.Range("C73").Value = -12.53
$ws.Range("C74").Value = -12.078
$ws.Range("D76").Value = -7.398000000000001
$ws.Range("C77").Value = -12.706
$ws.Range("D81").Value = -7.984999999999999
$ws.Range("D84").Value = -8.360000000000001
$ws.Range("D89").Value = -7.901000000000001
$ws.Range("C92").Value = -10.49
$ws.Range("D94").Value = -7.593999999999999
$ws.Range("C100").Value = -11.378

Write-Host "Applied 50 cell updates"
